# "Very rough calibration for Bulgaria"
# - constants!B2 (tb_n_contact): 15 -> 15.9
# - constants: delete row 7 (program_prop_child_reporting) entirely, shifting
#   rows 8:54 up to 7:53 (row 8's "age_breakpoints" data merges into row 7)
# - constants (after the shift) start_mdr_introduce_time: 1940 -> 1950
# - constants (after the shift) end_mdr_introduce_time:   1950 -> 1960
# - constants (after the shift) plot_start_time:           2010 -> 1990

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Bump the contact number parameter.
$ws.Range("B2").Value = 15.9

# Remove the "program_prop_child_reporting" row outright (row 7); this shifts
# every row below it up by one, matching the diff exactly (dimension
# A1:G54 -> A1:G53, shared string pruned, validations/ignoredErrors shift).
$ws.Rows("7:7").Delete()

# After the deletion the MDR-introduction window and plotting start year move.
$ws.Range("B14").Value = 1950
$ws.Range("B15").Value = 1960
$ws.Range("B16").Value = 1990
